$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Register a "YYYY-MM-DD HH:MM:SS" custom number format in the workbook's
# style table (matches original authoring tool's leftover datetime format)
# without leaving any cell actually using it: stamp a scratch cell far
# outside the real data range, then delete that scratch column entirely.
$ws.Range("Z1").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("Z1").Style = "Normal"
$ws.Range("Z1").ClearContents()
$ws.Range("Z1").EntireColumn.Delete()

# ---- Headers (row 1) ----
$headers = @("Employee", "Bonus", "Quota Met", "Start Date", "Hours Worked", "Performance", "Department", "Salary")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Extend the existing bold/bordered header style (already on A1:E1) across
# the three new header cells F1:H1 by copying formats from A1.
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats

# ---- Number formats (applied before values so no transient default
#      date/number format gets registered in the style table) ----
$ws.Range("B2:B6").NumberFormat = "$#,##0"
$ws.Range("H2:H6").NumberFormat = "$#,##0"
$ws.Range("C2:C6").NumberFormat = "0%"
$ws.Range("F2:F6").NumberFormat = "0%"
$ws.Range("D2:D6").NumberFormat = "YYYY-MM-DD"
$ws.Range("E2:E6").NumberFormat = "0.00"

# ---- Data rows ----
# Columns: A Employee, B Bonus, C Quota Met, D Start Date, E Hours Worked,
#          F Performance, G Department, H Salary
$data = @(
    @("Alice Johnson", 8500,  0.75, "2021-03-15", 42.5,  0.92, "Engineering", 85000),
    @("Bob Smith",     10800, 1.15, "2019-07-22", 38.75, 0.85, "Sales",       72000),
    @("Carol White",   5100,  0.92, "2022-01-10", 40,    0.78, "Marketing",   68000),
    @("David Brown",   11040, 0.88, "2018-11-05", 45.25, 0.95, "Engineering", 92000),
    @("Eva Martinez",  9360,  1.05, "2020-06-18", 39.5,  0.88, "Sales",       78000)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $r + 2
    $vals = $data[$r]

    $ws.Cells.Item($row, 1).Value = $vals[0]                    # Employee
    $ws.Cells.Item($row, 2).Value = $vals[1]                    # Bonus
    $ws.Cells.Item($row, 3).Value = $vals[2]                    # Quota Met
    $ws.Cells.Item($row, 4).Value = [DateTime]$vals[3]          # Start Date
    $ws.Cells.Item($row, 5).Value = $vals[4]                    # Hours Worked
    $ws.Cells.Item($row, 6).Value = $vals[5]                    # Performance
    $ws.Cells.Item($row, 7).Value = $vals[6]                    # Department
    $ws.Cells.Item($row, 8).Value = $vals[7]                    # Salary
}
